# OTCRecommendations.xlsx edit script
# - Rename "Pain Control" sheet to "Pain_Control"
# - Clear out the unused "Question N" placeholder labels on the GERD sheet (A6:A13)
#   and on the Allergies sheet (A13:A19)
# - Fill in the missing "No" answers (and one "Yes") on the Allergies sheet
# - Re-point selections / the active tab to match the new state of the workbook

$wb = $excel.ActiveWorkbook

$wsGERD        = $wb.Worksheets.Item("GERD")
$wsAllergies   = $wb.Worksheets.Item("Allergies")
$wsPainControl = $wb.Worksheets.Item("Pain Control")
$wsConstipation = $wb.Worksheets.Item("Constipation")

# --- Rename the "Pain Control" sheet ---
$wsPainControl.Name = "Pain_Control"

# --- GERD sheet: clear the leftover "Question 5".."Question 12" placeholders ---
$wsGERD.Range("A6:A13").ClearContents()

# --- Allergies sheet: fill in the missing Yes/No answers ---
$wsAllergies.Range("C2:C11").Value = "No"
$wsAllergies.Range("B12").Value = "Yes"
$wsAllergies.Range("C12").Value = "No"

# Match the formatting already used by the other "No" cells in this column
# (cell F3 already carries that style and is otherwise empty/unused).
$wsAllergies.Range("F3").Copy()
$wsAllergies.Range("C2:C12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Clear the leftover "Question 6".."Question 12" placeholders
$wsAllergies.Range("A13:A19").ClearContents()

# --- Update selections on each sheet ---
$wsGERD.Range("C19").Select()
$wsPainControl.Range("E36").Select()
$wsConstipation.Range("D24").Select()

# Allergies becomes the active sheet, with B36 selected
$wsAllergies.Activate()
$wsAllergies.Range("B36").Select()
